$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1186.8148
$ws.Range("I28").Value = 709.1053000000001
$ws.Range("J28").Value = 2321.375
$ws.Range("K28").Value = 709.1053000000001
$ws.Range("L28").Value = 2321.375
$ws.Range("M28").Value = -224.1053000000001
$ws.Range("N28").Value = -3291.375

$ws.Range("H75").Value = 44900
$ws.Range("J75").Value = 44900
$ws.Range("L75").Value = 44900
$ws.Range("N75").Value = -46772

$ws.Range("H78").Value = 44900
$ws.Range("J78").Value = 44900
$ws.Range("L78").Value = 134700
$ws.Range("N78").Value = -144060

$ws.Range("H88").Value = 1155.0834
$ws.Range("I88").Value = 1091.4
$ws.Range("J88").Value = 1200.5714
$ws.Range("K88").Value = 1091.4
$ws.Range("L88").Value = 1200.5714
$ws.Range("M88").Value = -685.4000000000001
$ws.Range("N88").Value = -2012.5714

$ws.Range("H91").Value = 1155.0834
$ws.Range("I91").Value = 1091.4
$ws.Range("J91").Value = 1200.5714
$ws.Range("K91").Value = 1091.4
$ws.Range("L91").Value = 1200.5714
$ws.Range("M91").Value = 312.5999999999999
$ws.Range("N91").Value = -4008.5714

$ws.Range("H96").Value = 1085.8
$ws.Range("I96").Value = 367
$ws.Range("J96").Value = 2164
$ws.Range("K96").Value = 1101
$ws.Range("L96").Value = 6492
$ws.Range("M96").Value = 272
$ws.Range("N96").Value = -9238

$ws.Range("H112").Value = 3646.6667
$ws.Range("I112").Value = 1090
$ws.Range("J112").Value = 4040
$ws.Range("K112").Value = 3270
$ws.Range("L112").Value = 12120
$ws.Range("M112").Value = -2162
$ws.Range("N112").Value = -14336

$ws.Range("H125").Value = 3784.1428
$ws.Range("J125").Value = 3373.25
$ws.Range("L125").Value = 30359.25
$ws.Range("N125").Value = -35279.25

$ws.Range("H131").Value = 69133
$ws.Range("I131").Value = 69133
$ws.Range("K131").Value = 207399
$ws.Range("M131").Value = -202359

$ws.Range("H132").Value = 2737.4666
$ws.Range("I132").Value = 2839.244
$ws.Range("J132").Value = 1694.25
$ws.Range("K132").Value = 8517.732
$ws.Range("L132").Value = 5082.75
$ws.Range("M132").Value = -5987.732
$ws.Range("N132").Value = -10142.75

$ws.Range("H135").Value = 7458.737
$ws.Range("I135").Value = 2317.611
$ws.Range("J135").Value = 99999
$ws.Range("K135").Value = 20858.499
$ws.Range("L135").Value = 899991
$ws.Range("M135").Value = -18323.499
$ws.Range("N135").Value = -905061

$ws.Range("H138").Value = 4837.0894
$ws.Range("I138").Value = 1006.75
$ws.Range("J138").Value = 5475.479
$ws.Range("K138").Value = 3020.25
$ws.Range("L138").Value = 16426.437
$ws.Range("M138").Value = 2119.75
$ws.Range("N138").Value = -26706.437

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6260.6064
$ws.Range("I32").Value = 662.9375
$ws.Range("J32").Value = 26928.924
$ws.Range("K32").Value = 662.9375
$ws.Range("L32").Value = 26928.924
$ws.Range("M32").Value = -375.9375
$ws.Range("N32").Value = -27502.924

$ws.Range("H61").Value = 4316
$ws.Range("I61").Value = 4337.72
$ws.Range("J61").Value = 4261.7
$ws.Range("K61").Value = 4337.72
$ws.Range("L61").Value = 4261.7
$ws.Range("M61").Value = -4125.72
$ws.Range("N61").Value = -4685.7

$ws.Range("H88").Value = 2543.7778
$ws.Range("I88").Value = 2132.4285
$ws.Range("J88").Value = 3983.5
$ws.Range("K88").Value = 2132.4285
$ws.Range("L88").Value = 3983.5
$ws.Range("M88").Value = -1726.4285
$ws.Range("N88").Value = -4795.5

$ws.Range("H91").Value = 2543.7778
$ws.Range("I91").Value = 2132.4285
$ws.Range("J91").Value = 3983.5
$ws.Range("K91").Value = 2132.4285
$ws.Range("L91").Value = 3983.5
$ws.Range("M91").Value = -728.4285
$ws.Range("N91").Value = -6791.5

$ws.Range("H136").Value = 4316
$ws.Range("I136").Value = 4337.72
$ws.Range("J136").Value = 4261.7
$ws.Range("K136").Value = 13013.16
$ws.Range("L136").Value = 12785.1
$ws.Range("M136").Value = -10463.16
$ws.Range("N136").Value = -17885.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2742.9546
$ws.Range("I20").Value = 3162.5454
$ws.Range("J20").Value = 2323.3635
$ws.Range("K20").Value = 3162.5454
$ws.Range("L20").Value = 2323.3635
$ws.Range("M20").Value = -2915.5454
$ws.Range("N20").Value = -2817.3635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7018.778
$ws.Range("I31").Value = 9024.799999999999
$ws.Range("J31").Value = 4511.25
$ws.Range("K31").Value = 9024.799999999999
$ws.Range("L31").Value = 4511.25
$ws.Range("M31").Value = -8729.799999999999
$ws.Range("N31").Value = -5101.25

$ws.Range("H34").Value = 7018.778
$ws.Range("I34").Value = 9024.799999999999
$ws.Range("J34").Value = 4511.25
$ws.Range("K34").Value = 9024.799999999999
$ws.Range("L34").Value = 4511.25
$ws.Range("M34").Value = -8822.799999999999
$ws.Range("N34").Value = -4915.25

$ws.Range("H58").Value = 5348.636
$ws.Range("I58").Value = 5910.4287
$ws.Range("J58").Value = 4365.5
$ws.Range("K58").Value = 5910.4287
$ws.Range("L58").Value = 4365.5
$ws.Range("M58").Value = -5707.4287
$ws.Range("N58").Value = -4771.5

$ws.Range("H105").Value = 1775.5
$ws.Range("I105").Value = 1700.6666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1700.6666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 46.33339999999998
$ws.Range("N105").Value = -5494

$ws.Range("H132").Value = 3900.8823
$ws.Range("I132").Value = 3107.6667
$ws.Range("J132").Value = 9850
$ws.Range("K132").Value = 9323.000100000001
$ws.Range("L132").Value = 29550
$ws.Range("M132").Value = -6793.000100000001
$ws.Range("N132").Value = -34610

$ws.Range("H136").Value = 5348.636
$ws.Range("I136").Value = 5910.4287
$ws.Range("J136").Value = 4365.5
$ws.Range("K136").Value = 17731.2861
$ws.Range("L136").Value = 13096.5
$ws.Range("M136").Value = -15181.2861
$ws.Range("N136").Value = -18196.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1149.6
$ws.Range("I5").Value = 1154.5714
$ws.Range("K5").Value = 3463.7142
$ws.Range("M5").Value = -3351.7142
$ws.Range("N5").ClearContents()

$ws.Range("H14").Value = 4519.1904
$ws.Range("I14").Value = 4519.1904
$ws.Range("K14").Value = 13557.5712
$ws.Range("M14").Value = -13384.5712

$ws.Range("H32").Value = 1720.8334
$ws.Range("I32").Value = 1720.8334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5162.5002
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4879.5002
$ws.Range("N32").ClearContents()

$ws.Range("H41").Value = 1498
$ws.Range("I41").Value = 1498
$ws.Range("K41").Value = 4494
$ws.Range("M41").Value = -4156

$ws.Range("H42").Value = 3500
$ws.Range("J42").Value = 3500
$ws.Range("L42").Value = 10500
$ws.Range("N42").Value = -11568

$ws.Range("H46").Value = 84378.75
$ws.Range("I46").Value = 541
$ws.Range("J46").Value = 112324.664
$ws.Range("K46").Value = 1623
$ws.Range("L46").Value = 336973.992
$ws.Range("M46").Value = -1532
$ws.Range("N46").Value = -337155.992

$ws.Range("H70").Value = 12282.375
$ws.Range("I70").Value = 8652.200000000001
$ws.Range("K70").Value = 25956.6
$ws.Range("M70").Value = -25641.6
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 12282.375
$ws.Range("I73").Value = 8652.200000000001
$ws.Range("K73").Value = 25956.6
$ws.Range("M73").Value = -24864.6
$ws.Range("N73").ClearContents()

$ws.Range("H82").Value = 2999
$ws.Range("I82").Value = 2999
$ws.Range("K82").Value = 8997
$ws.Range("M82").Value = -8591

$ws.Range("H85").Value = 2999
$ws.Range("I85").Value = 2999
$ws.Range("K85").Value = 8997
$ws.Range("M85").Value = -7593

$ws.Range("H107").Value = 815.4783
$ws.Range("I107").Value = 616.375
$ws.Range("J107").Value = 921.6667
$ws.Range("K107").Value = 1849.125
$ws.Range("L107").Value = 2765.0001
$ws.Range("M107").Value = 70.875
$ws.Range("N107").Value = -6605.0001

$ws.Range("H135").Value = 1149.6
$ws.Range("I135").Value = 1154.5714
$ws.Range("K135").Value = 10391.1426
$ws.Range("M135").Value = -7856.142600000001
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 5500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5230
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 5500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4564
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 2952.0952
$ws.Range("I132").Value = 2644.375
$ws.Range("J132").Value = 3936.8
$ws.Range("K132").Value = 7933.125
$ws.Range("L132").Value = 11810.4
$ws.Range("M132").Value = -5403.125
$ws.Range("N132").Value = -16870.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6887.143
$ws.Range("I136").Value = 6887.1665
$ws.Range("J136").Value = 6887
$ws.Range("K136").Value = 20661.4995
$ws.Range("L136").Value = 20661
$ws.Range("M136").Value = -18111.4995
$ws.Range("N136").Value = -25761

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2113.125
$ws.Range("I132").Value = 1523.9231
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 4571.7693
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -2041.7693
$ws.Range("N132").Value = -19059.0005

$ws.Range("H136").Value = 4653.36
$ws.Range("I136").Value = 4947.0454
$ws.Range("J136").Value = 2499.6667
$ws.Range("K136").Value = 14841.1362
$ws.Range("L136").Value = 7499.000100000001
$ws.Range("M136").Value = -12291.1362
$ws.Range("N136").Value = -12599.0001
